# Apply updated market price / profit figures to the Leve profit sheets
# (values refreshed by the scheduled market-data runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2597797.2
$ws.Range("I6").Value = 2857477
$ws.Range("K6").Value = 8572431
$ws.Range("M6").Value = -8572319
$ws.Range("H33").Value = 643.34485
$ws.Range("I33").Value = 663.8214
$ws.Range("K33").Value = 663.8214
$ws.Range("M33").Value = -434.8214
$ws.Range("H53").Value = 485.66666
$ws.Range("I53").Value = 209.21428
$ws.Range("J53").Value = 1038.5714
$ws.Range("K53").Value = 209.21428
$ws.Range("L53").Value = 1038.5714
$ws.Range("M53").Value = 427.78572
$ws.Range("N53").Value = -2312.5714
$ws.Range("H74").Value = 4337.9287
$ws.Range("I74").Value = 4415.0835
$ws.Range("J74").Value = 3875
$ws.Range("K74").Value = 4415.0835
$ws.Range("L74").Value = 3875
$ws.Range("M74").Value = -3479.0835
$ws.Range("N74").Value = -5747
$ws.Range("H77").Value = 4337.9287
$ws.Range("I77").Value = 4415.0835
$ws.Range("J77").Value = 3875
$ws.Range("K77").Value = 22075.4175
$ws.Range("L77").Value = 19375
$ws.Range("M77").Value = -17395.4175
$ws.Range("N77").Value = -28735
$ws.Range("H101").Value = 2233.75
$ws.Range("I101").Value = 2160
$ws.Range("J101").Value = 2267.2727
$ws.Range("K101").Value = 6480
$ws.Range("L101").Value = 6801.8181
$ws.Range("M101").Value = -4858
$ws.Range("N101").Value = -10045.8181
$ws.Range("H120").Value = 36000
$ws.Range("J120").Value = 36000
$ws.Range("L120").Value = 36000
$ws.Range("N120").Value = -45676
$ws.Range("H132").Value = 5091.986
$ws.Range("I132").Value = 4354.585
$ws.Range("K132").Value = 13063.755
$ws.Range("M132").Value = -10533.755
$ws.Range("H135").Value = 922
$ws.Range("I135").Value = 875.4828
$ws.Range("J135").Value = 1090.625
$ws.Range("K135").Value = 7879.3452
$ws.Range("L135").Value = 9815.625
$ws.Range("M135").Value = -5344.3452
$ws.Range("N135").Value = -14885.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3465.1853
$ws.Range("I74").Value = 4220
$ws.Range("J74").Value = 1672.5
$ws.Range("K74").Value = 4220
$ws.Range("L74").Value = 1672.5
$ws.Range("M74").Value = -3346
$ws.Range("N74").Value = -3420.5
$ws.Range("H77").Value = 3465.1853
$ws.Range("I77").Value = 4220
$ws.Range("J77").Value = 1672.5
$ws.Range("K77").Value = 21100
$ws.Range("L77").Value = 8362.5
$ws.Range("M77").Value = -16732
$ws.Range("N77").Value = -17098.5
$ws.Range("H104").Value = 46245
$ws.Range("J104").Value = 46245
$ws.Range("L104").Value = 46245
$ws.Range("N104").Value = -53233
$ws.Range("H110").Value = 1991.75
$ws.Range("I110").Value = 1765.7693
$ws.Range("J110").Value = 2971
$ws.Range("K110").Value = 1765.7693
$ws.Range("L110").Value = 2971
$ws.Range("M110").Value = 279.2307000000001
$ws.Range("N110").Value = -7061
$ws.Range("H132").Value = 5338.814
$ws.Range("I132").Value = 3447.9
$ws.Range("J132").Value = 6983.087
$ws.Range("K132").Value = 10343.7
$ws.Range("L132").Value = 20949.261
$ws.Range("M132").Value = -7813.700000000001
$ws.Range("N132").Value = -26009.261

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 5044.7856
$ws.Range("I75").Value = 1162.7
$ws.Range("J75").Value = 14750
$ws.Range("K75").Value = 1162.7
$ws.Range("L75").Value = 14750
$ws.Range("M75").Value = -226.7
$ws.Range("N75").Value = -16622
$ws.Range("H78").Value = 5044.7856
$ws.Range("I78").Value = 1162.7
$ws.Range("J78").Value = 14750
$ws.Range("K78").Value = 3488.1
$ws.Range("L78").Value = 44250
$ws.Range("M78").Value = 1191.9
$ws.Range("N78").Value = -53610
$ws.Range("H105").Value = 2363.6
$ws.Range("I105").Value = 1920
$ws.Range("J105").Value = 3029
$ws.Range("K105").Value = 1920
$ws.Range("L105").Value = 3029
$ws.Range("M105").Value = -173
$ws.Range("N105").Value = -6523

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3902.6924
$ws.Range("I16").Value = 6280.5
$ws.Range("J16").Value = 2845.889
$ws.Range("K16").Value = 6280.5
$ws.Range("L16").Value = 2845.889
$ws.Range("M16").Value = -5993.5
$ws.Range("N16").Value = -3419.889
$ws.Range("H43").Value = 20922.8
$ws.Range("J43").Value = 20922.8
$ws.Range("L43").Value = 20922.8
$ws.Range("N43").Value = -21290.8
$ws.Range("H101").Value = 20922.8
$ws.Range("J101").Value = 20922.8
$ws.Range("L101").Value = 20922.8
$ws.Range("N101").Value = -27412.8
$ws.Range("H105").Value = 1218.0385
$ws.Range("I105").Value = 1257.4166
$ws.Range("J105").Value = 745.5
$ws.Range("K105").Value = 1257.4166
$ws.Range("L105").Value = 745.5
$ws.Range("M105").Value = 489.5834
$ws.Range("N105").Value = -4239.5
$ws.Range("H113").Value = 3902.6924
$ws.Range("I113").Value = 6280.5
$ws.Range("J113").Value = 2845.889
$ws.Range("K113").Value = 6280.5
$ws.Range("L113").Value = 2845.889
$ws.Range("M113").Value = -4110.5
$ws.Range("N113").Value = -7185.889
$ws.Range("H134").Value = 2242.087
$ws.Range("I134").Value = 1014.5
$ws.Range("J134").Value = 3581.2727
$ws.Range("K134").Value = 3043.5
$ws.Range("L134").Value = 10743.8181
$ws.Range("M134").Value = -508.5
$ws.Range("N134").Value = -15813.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 350
$ws.Range("I7").Value = 155
$ws.Range("J7").Value = 428
$ws.Range("K7").Value = 465
$ws.Range("L7").Value = 1284
$ws.Range("M7").Value = -353
$ws.Range("N7").Value = -1508
$ws.Range("H33").Value = 293.9091
$ws.Range("I33").Value = 95
$ws.Range("J33").Value = 459.66666
$ws.Range("K33").Value = 570
$ws.Range("L33").Value = 2757.99996
$ws.Range("M33").Value = -287
$ws.Range("N33").Value = -3323.99996
$ws.Range("H59").Value = 1997.6666
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1997.6666
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 5992.9998
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -7072.9998
$ws.Range("H80").Value = 2966.6667
$ws.Range("J80").Value = 2966.6667
$ws.Range("L80").Value = 8900.000100000001
$ws.Range("N80").Value = -10772.0001
$ws.Range("H83").Value = 2966.6667
$ws.Range("J83").Value = 2966.6667
$ws.Range("L83").Value = 26700.0003
$ws.Range("N83").Value = -36060.0003
$ws.Range("H92").Value = 1463.25
$ws.Range("I92").Value = 950
$ws.Range("J92").Value = 1634.3334
$ws.Range("K92").Value = 2850
$ws.Range("L92").Value = 4903.0002
$ws.Range("M92").Value = -1602
$ws.Range("N92").Value = -7399.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3140.72
$ws.Range("I7").Value = 2367.0908
$ws.Range("J7").Value = 3748.5715
$ws.Range("K7").Value = 2367.0908
$ws.Range("L7").Value = 3748.5715
$ws.Range("M7").Value = -2255.0908
$ws.Range("N7").Value = -3972.5715
$ws.Range("H94").Value = 15135.556
$ws.Range("J94").Value = 15135.556
$ws.Range("L94").Value = 15135.556
$ws.Range("N94").Value = -16487.556
$ws.Range("H126").Value = 3140.72
$ws.Range("I126").Value = 2367.0908
$ws.Range("J126").Value = 3748.5715
$ws.Range("K126").Value = 7101.2724
$ws.Range("L126").Value = 11245.7145
$ws.Range("M126").Value = -4631.2724
$ws.Range("N126").Value = -16185.7145

